$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "Advanced source sync (for Japan trip) - for when there is no ref photo shared by all sources, but only e.g. A-B + B-C + B-D + C-E"
$ws.Range("A19").Value = "`"View`" button for ref. photos so that user can quickly check if the right photo is selected (just open associated program is ok)"

$ws.Range("A20").Select()
